$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1799
$ws1.Range("F6").Value = 1166
$ws1.Range("F12").Value = 1713
$ws1.Range("F18").Value = 1629
$ws1.Range("F24").Value = 12449
$ws1.Range("F25").Value = 12501
$ws1.Range("F29").Value = 271
$ws1.Range("F32").Value = 1946
$ws1.Range("F33").Value = 10

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 1

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1799
$ws4.Range("F7").Value = 1166
$ws4.Range("F14").Value = 1713
$ws4.Range("F23").Value = 1629
$ws4.Range("F28").Value = 1
$ws4.Range("F30").Value = 12449
$ws4.Range("F31").Value = 12501
$ws4.Range("F35").Value = 271
$ws4.Range("F36").Value = 0
$ws4.Range("F40").Value = 1946
$ws4.Range("F41").Value = 10

$wb.Save()
